$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$oldGuid = "0c7a40ae-fa35-4e32-acf3-e8083352c534"
$newGuid = "6bbff5fe-2d23-4f51-baa8-e64361157706"

$oldZhHash = "07e40451adc1dfc56339231748313436b784e3f9"
$newZhHash = "49c5d8b27207b496a6b4397d0ff9e440d659de58"

$newFileName = $newGuid + ".md"
$newPathAndName = "e2e\" + $newGuid + ".md"
$newZhXlf = $newGuid + "." + $newZhHash + ".zh-cn.xlf"
$newDeXlf = $newGuid + "." + $newZhHash + ".de-de.xlf"

# --- Overview sheet ---
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Range("G2").Value = "2016-08-15 10:55:10"

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-15 10:54:58"

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-08-15 10:55:10"

# --- Hyperlink display text updates (keep same target, only the visible text changes) ---
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = $newPathAndName
}
foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = $newFileName
}
foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = $newFileName
}

Write-Host "Report regenerated for handoff"
